$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datasets")

# --- Update the "Block Group Housing (2013).Rds" row (row 5) ---
# B5: description text now mentions both block-group and county level.
$ws.Range("B5").Value = "Median year built and total number of housing units at both the Census block group and county level."

# C5: notes text gets an added second line (wrapped) explaining block groups.
$ws.Range("C5").Value = "Compiled from the 2013 vintage of the American Community Survey. `nNote that block groups generally contain 600–3,000 people, and never cross state, county, or census tract boundaries."

# Turn on word-wrap for the longer notes/description cells so the new text is readable.
$ws.Range("B5:C5").WrapText = $true

# Widen the Notes column (C) to better fit the longer note.
$ws.Columns.Item(3).ColumnWidth = 96

# Grow row 5 to fit the now-wrapped, two-line text.
$ws.Rows.Item(5).RowHeight = 23.85
